# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Leve profit tracking sheets
# (values refreshed by the scheduled market-data runner).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 76926540
$ws.Cells.Item(64, 9).Value = 1000000000
$ws.Cells.Item(64, 10).Value = 3755.0833
$ws.Cells.Item(64, 11).Value = 1000000000
$ws.Cells.Item(64, 12).Value = 3755.0833
$ws.Cells.Item(64, 13).Value = -999999752
$ws.Cells.Item(64, 14).Value = -4251.0833

$ws.Cells.Item(67, 8).Value = 76926540
$ws.Cells.Item(67, 9).Value = 1000000000
$ws.Cells.Item(67, 10).Value = 3755.0833
$ws.Cells.Item(67, 11).Value = 1000000000
$ws.Cells.Item(67, 12).Value = 3755.0833
$ws.Cells.Item(67, 13).Value = -999999142
$ws.Cells.Item(67, 14).Value = -5471.0833

$ws.Cells.Item(74, 8).Value = 3132.5
$ws.Cells.Item(74, 9).Value = 2338.8
$ws.Cells.Item(74, 10).Value = 3437.7693
$ws.Cells.Item(74, 11).Value = 2338.8
$ws.Cells.Item(74, 12).Value = 3437.7693
$ws.Cells.Item(74, 13).Value = -1402.8
$ws.Cells.Item(74, 14).Value = -5309.7693

$ws.Cells.Item(76, 8).Value = 5570.727
$ws.Cells.Item(76, 9).Value = 2880.6
$ws.Cells.Item(76, 10).Value = 7812.5
$ws.Cells.Item(76, 11).Value = 2880.6
$ws.Cells.Item(76, 12).Value = 7812.5
$ws.Cells.Item(76, 13).Value = -2565.6
$ws.Cells.Item(76, 14).Value = -8442.5

$ws.Cells.Item(77, 8).Value = 3132.5
$ws.Cells.Item(77, 9).Value = 2338.8
$ws.Cells.Item(77, 10).Value = 3437.7693
$ws.Cells.Item(77, 11).Value = 11694
$ws.Cells.Item(77, 12).Value = 17188.8465
$ws.Cells.Item(77, 13).Value = -7014
$ws.Cells.Item(77, 14).Value = -26548.8465

$ws.Cells.Item(79, 8).Value = 5570.727
$ws.Cells.Item(79, 9).Value = 2880.6
$ws.Cells.Item(79, 10).Value = 7812.5
$ws.Cells.Item(79, 11).Value = 2880.6
$ws.Cells.Item(79, 12).Value = 7812.5
$ws.Cells.Item(79, 13).Value = -1788.6
$ws.Cells.Item(79, 14).Value = -9996.5

$ws.Cells.Item(98, 8).Value = 1609.3939
$ws.Cells.Item(98, 9).Value = 1357.0333
$ws.Cells.Item(98, 10).Value = 4133
$ws.Cells.Item(98, 11).Value = 1357.0333
$ws.Cells.Item(98, 12).Value = 4133
$ws.Cells.Item(98, 13).Value = 140.9666999999999
$ws.Cells.Item(98, 14).Value = -7129

$ws.Cells.Item(116, 8).Value = 40579
$ws.Cells.Item(116, 9).Value = 6155
$ws.Cells.Item(116, 11).Value = 6155
$ws.Cells.Item(116, 13).Value = -2713

$ws.Cells.Item(122, 8).Value = 1609.3939
$ws.Cells.Item(122, 9).Value = 1357.0333
$ws.Cells.Item(122, 10).Value = 4133
$ws.Cells.Item(122, 11).Value = 4071.0999
$ws.Cells.Item(122, 12).Value = 12399
$ws.Cells.Item(122, 13).Value = -1621.0999
$ws.Cells.Item(122, 14).Value = -17299


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1720.15
$ws.Cells.Item(32, 9).Value = 1632.8062
$ws.Cells.Item(32, 10).Value = 6000
$ws.Cells.Item(32, 11).Value = 1632.8062
$ws.Cells.Item(32, 12).Value = 6000
$ws.Cells.Item(32, 13).Value = -1345.8062
$ws.Cells.Item(32, 14).Value = -6574

$ws.Cells.Item(63, 8).Value = 2213.125
$ws.Cells.Item(63, 9).Value = 2381
$ws.Cells.Item(63, 10).Value = 1933.3334
$ws.Cells.Item(63, 11).Value = 2381
$ws.Cells.Item(63, 12).Value = 1933.3334
$ws.Cells.Item(63, 13).Value = -1695
$ws.Cells.Item(63, 14).Value = -3305.3334

$ws.Cells.Item(66, 8).Value = 2213.125
$ws.Cells.Item(66, 9).Value = 2381
$ws.Cells.Item(66, 10).Value = 1933.3334
$ws.Cells.Item(66, 11).Value = 11905
$ws.Cells.Item(66, 12).Value = 9666.666999999999
$ws.Cells.Item(66, 13).Value = -8473
$ws.Cells.Item(66, 14).Value = -16530.667


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 34000
$ws.Cells.Item(26, 10).Value = 34000
$ws.Cells.Item(26, 12).Value = 34000
$ws.Cells.Item(26, 14).Value = -34584

$ws.Cells.Item(96, 8).Value = 18611.875
$ws.Cells.Item(96, 9).Value = 8105.25
$ws.Cells.Item(96, 10).Value = 29118.5
$ws.Cells.Item(96, 11).Value = 8105.25
$ws.Cells.Item(96, 12).Value = 29118.5
$ws.Cells.Item(96, 13).Value = -5359.25
$ws.Cells.Item(96, 14).Value = -34610.5

$ws.Cells.Item(105, 8).Value = 1714.08
$ws.Cells.Item(105, 9).Value = 1309.3529
$ws.Cells.Item(105, 10).Value = 2574.125
$ws.Cells.Item(105, 11).Value = 1309.3529
$ws.Cells.Item(105, 12).Value = 2574.125
$ws.Cells.Item(105, 13).Value = 437.6470999999999
$ws.Cells.Item(105, 14).Value = -6068.125

$ws.Cells.Item(112, 8).Value = 48000
$ws.Cells.Item(112, 10).Value = 48000
$ws.Cells.Item(112, 12).Value = 48000
$ws.Cells.Item(112, 14).Value = -50954


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1291.26
$ws.Cells.Item(31, 9).Value = 1052.0741
$ws.Cells.Item(31, 10).Value = 1572.0435
$ws.Cells.Item(31, 11).Value = 1052.0741
$ws.Cells.Item(31, 12).Value = 1572.0435
$ws.Cells.Item(31, 13).Value = -757.0741
$ws.Cells.Item(31, 14).Value = -2162.0435

$ws.Cells.Item(34, 8).Value = 1291.26
$ws.Cells.Item(34, 9).Value = 1052.0741
$ws.Cells.Item(34, 10).Value = 1572.0435
$ws.Cells.Item(34, 11).Value = 1052.0741
$ws.Cells.Item(34, 12).Value = 1572.0435
$ws.Cells.Item(34, 13).Value = -850.0741
$ws.Cells.Item(34, 14).Value = -1976.0435

$ws.Cells.Item(62, 8).Value = 4384.375
$ws.Cells.Item(62, 9).Value = 4812
$ws.Cells.Item(62, 10).Value = 2857.1428
$ws.Cells.Item(62, 11).Value = 4812
$ws.Cells.Item(62, 12).Value = 2857.1428
$ws.Cells.Item(62, 13).Value = -4188
$ws.Cells.Item(62, 14).Value = -4105.1428

$ws.Cells.Item(65, 8).Value = 4384.375
$ws.Cells.Item(65, 9).Value = 4812
$ws.Cells.Item(65, 10).Value = 2857.1428
$ws.Cells.Item(65, 11).Value = 24060
$ws.Cells.Item(65, 12).Value = 14285.714
$ws.Cells.Item(65, 13).Value = -20940
$ws.Cells.Item(65, 14).Value = -20525.714

$ws.Cells.Item(100, 8).Value = 24780
$ws.Cells.Item(100, 10).Value = 24780
$ws.Cells.Item(100, 12).Value = 24780
$ws.Cells.Item(100, 14).Value = -26944

$ws.Cells.Item(124, 8).Value = 18168.867
$ws.Cells.Item(124, 10).Value = 18168.867
$ws.Cells.Item(124, 12).Value = 18168.867
$ws.Cells.Item(124, 14).Value = -23078.867


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 177.3
$ws.Cells.Item(50, 9).Value = 131.6
$ws.Cells.Item(50, 10).Value = 223
$ws.Cells.Item(50, 11).Value = 394.8
$ws.Cells.Item(50, 12).Value = 669
$ws.Cells.Item(50, 13).Value = 86.20000000000005
$ws.Cells.Item(50, 14).Value = -1631

$ws.Cells.Item(53, 8).Value = 177.3
$ws.Cells.Item(53, 9).Value = 131.6
$ws.Cells.Item(53, 10).Value = 223
$ws.Cells.Item(53, 11).Value = 394.8
$ws.Cells.Item(53, 12).Value = 669
$ws.Cells.Item(53, 13).Value = 86.20000000000005
$ws.Cells.Item(53, 14).Value = -1631

$ws.Cells.Item(117, 8).Value = 716.6667
$ws.Cells.Item(117, 10).Value = 925
$ws.Cells.Item(117, 12).Value = 2775
$ws.Cells.Item(117, 14).Value = -9659

$ws.Cells.Item(121, 8).Value = 8793.924999999999
$ws.Cells.Item(121, 10).Value = 9969.343000000001
$ws.Cells.Item(121, 12).Value = 29908.029
$ws.Cells.Item(121, 14).Value = -32528.029

$ws.Cells.Item(137, 8).Value = 3004229.2
$ws.Cells.Item(137, 9).Value = 59162.473
$ws.Cells.Item(137, 10).Value = 6734647
$ws.Cells.Item(137, 11).Value = 177487.419
$ws.Cells.Item(137, 12).Value = 20203941
$ws.Cells.Item(137, 13).Value = -172387.419
$ws.Cells.Item(137, 14).Value = -20214141

$ws.Cells.Item(140, 8).Value = 26782.432
$ws.Cells.Item(140, 9).Value = 29828.59
$ws.Cells.Item(140, 10).Value = 3022.4
$ws.Cells.Item(140, 11).Value = 89485.77
$ws.Cells.Item(140, 12).Value = 9067.200000000001
$ws.Cells.Item(140, 13).Value = -84305.77
$ws.Cells.Item(140, 14).Value = -19427.2


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 3899.8572
$ws.Cells.Item(70, 9).Value = 3883.1667
$ws.Cells.Item(70, 10).Value = 4000
$ws.Cells.Item(70, 11).Value = 3883.1667
$ws.Cells.Item(70, 12).Value = 4000
$ws.Cells.Item(70, 13).Value = -3613.1667
$ws.Cells.Item(70, 14).Value = -4540

$ws.Cells.Item(73, 8).Value = 3899.8572
$ws.Cells.Item(73, 9).Value = 3883.1667
$ws.Cells.Item(73, 10).Value = 4000
$ws.Cells.Item(73, 11).Value = 3883.1667
$ws.Cells.Item(73, 12).Value = 4000
$ws.Cells.Item(73, 13).Value = -2947.1667
$ws.Cells.Item(73, 14).Value = -5872

$ws.Cells.Item(80, 8).Value = 2372.8572
$ws.Cells.Item(80, 9).Value = 2203.3333
$ws.Cells.Item(80, 10).Value = 2500
$ws.Cells.Item(80, 11).Value = 2203.3333
$ws.Cells.Item(80, 12).Value = 2500
$ws.Cells.Item(80, 13).Value = -1205.3333
$ws.Cells.Item(80, 14).Value = -4496

$ws.Cells.Item(83, 8).Value = 2372.8572
$ws.Cells.Item(83, 9).Value = 2203.3333
$ws.Cells.Item(83, 10).Value = 2500
$ws.Cells.Item(83, 11).Value = 11016.6665
$ws.Cells.Item(83, 12).Value = 12500
$ws.Cells.Item(83, 13).Value = -6024.666499999999
$ws.Cells.Item(83, 14).Value = -22484


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2440
$ws.Cells.Item(40, 9).Value = 1815
$ws.Cells.Item(40, 10).Value = 3377.5
$ws.Cells.Item(40, 11).Value = 1815
$ws.Cells.Item(40, 12).Value = 3377.5
$ws.Cells.Item(40, 13).Value = -1679
$ws.Cells.Item(40, 14).Value = -3649.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 8334211
$ws.Cells.Item(122, 9).Value = 8696490
$ws.Cells.Item(122, 11).Value = 26089470
$ws.Cells.Item(122, 13).Value = -26087020

$ws.Cells.Item(132, 8).Value = 1293.2924
$ws.Cells.Item(132, 9).Value = 1095.409
$ws.Cells.Item(132, 10).Value = 1707.9048
$ws.Cells.Item(132, 11).Value = 3286.227
$ws.Cells.Item(132, 12).Value = 5123.7144
$ws.Cells.Item(132, 13).Value = -756.2270000000003
$ws.Cells.Item(132, 14).Value = -10183.7144

